# Input code file name change
# - Adds a new "Graph" worksheet (a copy of the Python-code test-data table
#   that already lives on "DataStructure"), placed after "Login".
# - Updates the selection remembered on "DataStructure".
# - Leaves "Login" as a normal (non-active) tab since "Graph" becomes the
#   newly active / selected sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the remembered selection on the "DataStructure" sheet.
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("DataStructure")
$wsData.Activate()
$wsData.Range("A1:C3").Select()

# ---------------------------------------------------------------------
# 2. Add the new "Graph" worksheet right after "Login" and populate it
#    with the same TestCaseID / Python Code / Expected Output table.
# ---------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("Login")
$wsGraph = $wb.Worksheets.Add($null, $wsLogin)
$wsGraph.Name = "Graph"

$wsGraph.Range("A1").Value = "TestCaseID"
$wsGraph.Range("B1").Value = "Python Code"
$wsGraph.Range("C1").Value = "Expected Output"

$wsGraph.Range("A2").Value = "ValidCode"
$wsGraph.Range("B2").Value = 'print("hello");'
$wsGraph.Range("C2").Value = "hello"

$wsGraph.Range("A3").Value = "InvalidCode"
$wsGraph.Range("B3").Value = "xyz"
$wsGraph.Range("C3").Value = "xyz not defined on line 1"

# Size the columns to fit their contents (matches the "bestFit" columns
# used for the same data on "DataStructure").
$wsGraph.Columns.Item(1).AutoFit()
$wsGraph.Columns.Item(2).AutoFit()
$wsGraph.Columns.Item(3).AutoFit()

# Remembered selection on the new sheet.
$wsGraph.Range("C5").Select()
